# Update Fgf1-Cd44 LR-pair sheet with new TPM-based values (adds Resolving-Mac cluster)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ("ECs","Fgf1","Cd44","ECs",3,1,1.190640666666667,3.571922,0.1136540143525372,0.1136540143525372,3,1,24.576554,73.729662,0.07553767049546639,0.07553767049546638,29.26184463892934,263.3566017503641,0.008585159486648966,0.008585159486648964),
    ("ECs","Fgf1","Cd44","FAPs",3,1,1.190640666666667,3.571922,0.1136540143525372,0.1136540143525372,3,1,117.0512696666667,351.153809,0.359764849016532,0.359764849016532,139.3660017500998,1254.294015750898,0.04088871931366332,0.04088871931366332),
    ("ECs","Fgf1","Cd44","MuSCs",3,1,1.190640666666667,3.571922,0.1136540143525372,0.1136540143525372,3,1,55.68784966666667,167.063549,0.1711603033819035,0.1711603033819035,66.30421845235311,596.7379660711781,0.01945305557715149,0.01945305557715149),
    ("ECs","Fgf1","Cd44","Resolving-Mac",3,1,1.190640666666667,3.571922,0.1136540143525372,0.1136540143525372,3,1,128.0392633333333,384.11779,0.3935371771060981,0.3935371771060981,152.4487538547089,1372.03878469238,0.04472707997507346,0.04472707997507346),
    ("FAPs","Fgf1","Cd44","ECs",3,1,5.168173666666667,15.504521,0.4933341355895272,0.4933341355895272,3,1,24.576554,73.729662,0.07553767049546639,0.07553767049546638,127.0158992002114,1143.143092801902,0.03726531137832745,0.03726531137832744),
    ("FAPs","Fgf1","Cd44","FAPs",3,1,5.168173666666667,15.504521,0.4933341355895272,0.4933341355895272,3,1,117.0512696666667,351.153809,0.359764849016532,0.359764849016532,604.9412895411655,5444.47160587049,0.1774842808050676,0.1774842808050676),
    ("FAPs","Fgf1","Cd44","MuSCs",3,1,5.168173666666667,15.504521,0.4933341355895272,0.4933341355895272,3,1,55.68784966666667,167.063549,0.1711603033819035,0.1711603033819035,287.8044782005588,2590.240303805029,0.0844392203161526,0.0844392203161526),
    ("FAPs","Fgf1","Cd44","Resolving-Mac",3,1,5.168173666666667,15.504521,0.4933341355895272,0.4933341355895272,3,1,128.0392633333333,384.11779,0.3935371771060981,0.3935371771060981,661.7291490587323,5955.56234152859,0.1941453230899796,0.1941453230899796),
    ("MuSCs","Fgf1","Cd44","ECs",3,1,4.037194,12.111582,0.3853751326204581,0.3853751326204581,3,1,24.576554,73.729662,0.07553767049546639,0.07553767049546638,99.22031634947602,892.982847145284,0.02911033978503082,0.02911033978503082),
    ("MuSCs","Fgf1","Cd44","FAPs",3,1,4.037194,12.111582,0.3853751326204581,0.3853751326204581,3,1,117.0512696666667,351.153809,0.359764849016532,0.359764849016532,472.5586835906487,4253.028152315838,0.1386444264019251,0.1386444264019251),
    ("MuSCs","Fgf1","Cd44","MuSCs",3,1,4.037194,12.111582,0.3853751326204581,0.3853751326204581,3,1,55.68784966666667,167.063549,0.1711603033819035,0.1711603033819035,224.8226525471687,2023.403872924518,0.06596092461515891,0.06596092461515891),
    ("MuSCs","Fgf1","Cd44","Resolving-Mac",3,1,4.037194,12.111582,0.3853751326204581,0.3853751326204581,3,1,128.0392633333333,384.11779,0.3935371771060981,0.3935371771060981,516.9193456937534,4652.274111243781,0.1516594418183433,0.1516594418183433),
    ("Resolving-Mac","Fgf1","Cd44","ECs",1,0.3333333333333333,0.08000233333333333,0.240007,0.007636717437477471,0.007636717437477472,3,1,24.576554,73.729662,0.07553767049546639,0.07553767049546638,1.966181665292667,17.695634987634,0.0005768598454591557,0.0005768598454591557),
    ("Resolving-Mac","Fgf1","Cd44","FAPs",1,0.3333333333333333,0.08000233333333333,0.240007,0.007636717437477471,0.007636717437477472,3,1,117.0512696666667,351.153809,0.359764849016532,0.359764849016532,9.364374692962555,84.27937223666301,0.002747422495875999,0.002747422495876),
    ("Resolving-Mac","Fgf1","Cd44","MuSCs",1,0.3333333333333333,0.08000233333333333,0.240007,0.007636717437477471,0.007636717437477472,3,1,55.68784966666667,167.063549,0.1711603033819035,0.1711603033819035,4.455157911649223,40.096421204843,0.001307102873440517,0.001307102873440517),
    ("Resolving-Mac","Fgf1","Cd44","Resolving-Mac",1,0.3333333333333333,0.08000233333333333,0.240007,0.007636717437477471,0.007636717437477472,3,1,128.0392633333333,384.11779,0.3935371771060981,0.3935371771060981,10.24343982494778,92.19095842453001,0.003005332222701799,0.0030053322227018) 
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

